# Auto-generated edit script applying the Valefor_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# for specific rows across multiple class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(93, 8).Value = 29601
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 29601
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 29601
$ws.Cells.Item(93, 14).Value = -34593
$ws.Cells.Item(132, 8).Value = 1737717.1
$ws.Cells.Item(132, 9).Value = 2137965.5
$ws.Cells.Item(132, 10).Value = 3307.6667
$ws.Cells.Item(132, 11).Value = 6413896.5
$ws.Cells.Item(132, 12).Value = 9923.000100000001
$ws.Cells.Item(132, 13).Value = -6411366.5
$ws.Cells.Item(132, 14).Value = -14983.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1830.0588
$ws.Cells.Item(61, 9).Value = 1716.2307
$ws.Cells.Item(61, 10).Value = 2200
$ws.Cells.Item(61, 11).Value = 1716.2307
$ws.Cells.Item(61, 12).Value = 2200
$ws.Cells.Item(61, 13).Value = -1504.2307
$ws.Cells.Item(61, 14).Value = -2624
$ws.Cells.Item(122, 8).Value = 1585.963
$ws.Cells.Item(122, 9).Value = 1531.5769
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 4594.7307
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -2144.7307
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(132, 8).Value = 1578.4546
$ws.Cells.Item(132, 9).Value = 858.125
$ws.Cells.Item(132, 10).Value = 3499.3333
$ws.Cells.Item(132, 11).Value = 2574.375
$ws.Cells.Item(132, 12).Value = 10497.9999
$ws.Cells.Item(132, 13).Value = -44.375
$ws.Cells.Item(132, 14).Value = -15557.9999
$ws.Cells.Item(136, 8).Value = 1830.0588
$ws.Cells.Item(136, 9).Value = 1716.2307
$ws.Cells.Item(136, 10).Value = 2200
$ws.Cells.Item(136, 11).Value = 5148.6921
$ws.Cells.Item(136, 12).Value = 6600
$ws.Cells.Item(136, 13).Value = -2598.6921
$ws.Cells.Item(136, 14).Value = -11700

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 366.66666
$ws.Cells.Item(22, 9).Value = 366.66666
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 366.66666
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -193.66666
$ws.Cells.Item(105, 8).Value = 3789395.2
$ws.Cells.Item(105, 9).Value = 4546674.5
$ws.Cells.Item(105, 10).Value = 3000
$ws.Cells.Item(105, 11).Value = 4546674.5
$ws.Cells.Item(105, 12).Value = 3000
$ws.Cells.Item(105, 13).Value = -4544927.5
$ws.Cells.Item(105, 14).Value = -6494
$ws.Cells.Item(134, 8).Value = 1754.9
$ws.Cells.Item(134, 9).Value = 1642.3429
$ws.Cells.Item(134, 10).Value = 2542.8
$ws.Cells.Item(134, 11).Value = 4927.028700000001
$ws.Cells.Item(134, 12).Value = 7628.400000000001
$ws.Cells.Item(134, 13).Value = -2392.028700000001
$ws.Cells.Item(134, 14).Value = -12698.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 62500692
$ws.Cells.Item(22, 9).Value = 71429224
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 71429224
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -71428874
$ws.Cells.Item(22, 14).Value = -1700
$ws.Cells.Item(88, 8).Value = 36307.168
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 36307.168
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 36307.168
$ws.Cells.Item(88, 14).Value = -37119.168
$ws.Cells.Item(91, 8).Value = 36307.168
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 36307.168
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 36307.168
$ws.Cells.Item(91, 14).Value = -39115.168
$ws.Cells.Item(132, 8).Value = 1987.3334
$ws.Cells.Item(132, 9).Value = 1194.6
$ws.Cells.Item(132, 10).Value = 2978.25
$ws.Cells.Item(132, 11).Value = 3583.8
$ws.Cells.Item(132, 12).Value = 8934.75
$ws.Cells.Item(132, 13).Value = -1053.8
$ws.Cells.Item(132, 14).Value = -13994.75
$ws.Cells.Item(134, 8).Value = 1812.7059
$ws.Cells.Item(134, 9).Value = 1440.2051
$ws.Cells.Item(134, 10).Value = 3023.3333
$ws.Cells.Item(134, 11).Value = 4320.615299999999
$ws.Cells.Item(134, 12).Value = 9069.999899999999
$ws.Cells.Item(134, 13).Value = -1785.615299999999
$ws.Cells.Item(134, 14).Value = -14139.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 5997.5
$ws.Cells.Item(94, 9).Value = 2000
$ws.Cells.Item(94, 10).Value = 7330
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 21990
$ws.Cells.Item(94, 13).Value = -5324
$ws.Cells.Item(94, 14).Value = -23342
$ws.Cells.Item(99, 8).Value = 2055.625
$ws.Cells.Item(99, 9).Value = 1590.8334
$ws.Cells.Item(99, 10).Value = 3450
$ws.Cells.Item(99, 11).Value = 4772.5002
$ws.Cells.Item(99, 12).Value = 10350
$ws.Cells.Item(99, 13).Value = -2526.5002
$ws.Cells.Item(99, 14).Value = -14842
$ws.Cells.Item(109, 8).Value = 100
$ws.Cells.Item(109, 9).Value = 100
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 300
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = 740
$ws.Cells.Item(112, 8).Value = 3067.6843
$ws.Cells.Item(112, 9).Value = 1425
$ws.Cells.Item(112, 10).Value = 3505.7334
$ws.Cells.Item(112, 11).Value = 4275
$ws.Cells.Item(112, 12).Value = 10517.2002
$ws.Cells.Item(112, 13).Value = -3167
$ws.Cells.Item(112, 14).Value = -12733.2002
$ws.Cells.Item(116, 8).Value = 1642.8572
$ws.Cells.Item(116, 9).Value = 875
$ws.Cells.Item(116, 10).Value = 2666.6667
$ws.Cells.Item(116, 11).Value = 2625
$ws.Cells.Item(116, 12).Value = 8000.000100000001
$ws.Cells.Item(116, 13).Value = 817
$ws.Cells.Item(116, 14).Value = -14884.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3429.4546
$ws.Cells.Item(132, 9).Value = 3290.625
$ws.Cells.Item(132, 10).Value = 3799.6667
$ws.Cells.Item(132, 11).Value = 9871.875
$ws.Cells.Item(132, 12).Value = 11399.0001
$ws.Cells.Item(132, 13).Value = -7341.875
$ws.Cells.Item(132, 14).Value = -16459.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1855.2413
$ws.Cells.Item(7, 9).Value = 1872.08
$ws.Cells.Item(7, 10).Value = 1750
$ws.Cells.Item(7, 11).Value = 1872.08
$ws.Cells.Item(7, 12).Value = 1750
$ws.Cells.Item(7, 13).Value = -1760.08
$ws.Cells.Item(7, 14).Value = -1974
$ws.Cells.Item(22, 8).Value = 374.15152
$ws.Cells.Item(22, 9).Value = 363.5926
$ws.Cells.Item(22, 10).Value = 421.66666
$ws.Cells.Item(22, 11).Value = 363.5926
$ws.Cells.Item(22, 12).Value = 421.66666
$ws.Cells.Item(22, 13).Value = -68.5926
$ws.Cells.Item(22, 14).Value = -1011.66666
$ws.Cells.Item(27, 8).Value = 374.15152
$ws.Cells.Item(27, 9).Value = 363.5926
$ws.Cells.Item(27, 10).Value = 421.66666
$ws.Cells.Item(27, 11).Value = 363.5926
$ws.Cells.Item(27, 12).Value = 421.66666
$ws.Cells.Item(27, 13).Value = -256.5926
$ws.Cells.Item(27, 14).Value = -635.66666
$ws.Cells.Item(111, 8).Value = 46193.5
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 46193.5
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 46193.5
$ws.Cells.Item(111, 14).Value = -54373.5
$ws.Cells.Item(126, 8).Value = 1855.2413
$ws.Cells.Item(126, 9).Value = 1872.08
$ws.Cells.Item(126, 10).Value = 1750
$ws.Cells.Item(126, 11).Value = 5616.24
$ws.Cells.Item(126, 12).Value = 5250
$ws.Cells.Item(126, 13).Value = -3146.24
$ws.Cells.Item(126, 14).Value = -10190
$ws.Cells.Item(132, 8).Value = 2842.158
$ws.Cells.Item(132, 9).Value = 2187.6875
$ws.Cells.Item(132, 10).Value = 6332.6665
$ws.Cells.Item(132, 11).Value = 6563.0625
$ws.Cells.Item(132, 12).Value = 18997.9995
$ws.Cells.Item(132, 13).Value = -4033.0625
$ws.Cells.Item(132, 14).Value = -24057.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 37300
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 37300
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 37300
$ws.Cells.Item(95, 13).ClearContents()
$ws.Cells.Item(95, 14).Value = -42792
$ws.Cells.Item(96, 8).Value = 1568.25
$ws.Cells.Item(96, 9).Value = 847.6
$ws.Cells.Item(96, 10).Value = 2769.3333
$ws.Cells.Item(96, 11).Value = 847.6
$ws.Cells.Item(96, 12).Value = 2769.3333
$ws.Cells.Item(96, 13).Value = 525.4
$ws.Cells.Item(96, 14).Value = -5515.3333
$ws.Cells.Item(122, 8).Value = 2078.6667
$ws.Cells.Item(122, 9).Value = 1252.5333
$ws.Cells.Item(122, 10).Value = 3455.5557
$ws.Cells.Item(122, 11).Value = 3757.5999
$ws.Cells.Item(122, 12).Value = 10366.6671
$ws.Cells.Item(122, 13).Value = -1307.5999
$ws.Cells.Item(122, 14).Value = -15266.6671
$ws.Cells.Item(126, 8).Value = 967.5294
$ws.Cells.Item(126, 9).Value = 674.3
$ws.Cells.Item(126, 10).Value = 1386.4286
$ws.Cells.Item(126, 11).Value = 2022.9
$ws.Cells.Item(126, 12).Value = 4159.2858
$ws.Cells.Item(126, 13).Value = 447.1000000000001
$ws.Cells.Item(126, 14).Value = -9099.2858
